$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liquidación")

$ws.Range("C14").Value = (Get-Date -Year 2023 -Month 6 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J14").Value = 2500000
$ws.Range("J15").Value = 150000
$ws.Range("C18").Value = 6
$ws.Range("C21").Value = (Get-Date -Year 2023 -Month 1 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C22").Value = (Get-Date -Year 2023 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C26").Value = (Get-Date -Year 2023 -Month 1 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C27").Value = (Get-Date -Year 2009 -Month 6 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C29").Value = 195
